# Generate Report for Handback
#
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - The zh-cn and de-de sheets get their "Latest Target File" (I) / "Latest
#    Handback File" (J) / "Latest Handback DateTime" (K) columns populated,
#    with the new Target File column becoming a hyperlink (like column A).
#  - A few report columns are widened so the new hyperlink text fits.

$wb = $excel.ActiveWorkbook

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e76aea10b0f5b2f081db216bba4ac45b5c5c420/e2e/24966354-00aa-4c6f-b8f1-c3a0a34cee9a.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7e76aea10b0f5b2f081db216bba4ac45b5c5c420/e2e/e5a11050-9d9c-4487-a705-04263cfa461a.md"

$newStatus = "Handed back: in sync with en-US"

# Column width that lands as close as possible to the report's target
# character width (the host quantizes stored width to input + 5/6).
$wideColWidth = 29.15
$fullColWidth = 235 / 6   # -> stored width of exactly 40

# ---------------------------------------------------------------------
# 1. Overview sheet: status text + widen the per-locale status columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("J2").Value = "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.252e131be62cc623e431e38256b7e8c7a708a19d.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-22 10:42:20"

$wsZh.Range("J3").Value = "e5a11050-9d9c-4487-a705-04263cfa461a.e97206c15f67c039fb036d86b15dd20e106ba7c3.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-22 10:42:20"

# Rebuild hyperlinks so the new "Latest Target File" links (column I) are
# interleaved with the existing ones in the same order Excel would assign
# relationship ids: A2, I2, A3, I3.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlMd1, [Type]::Missing, "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.md", "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlMd1, [Type]::Missing, "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.md", "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlMd2, [Type]::Missing, "e5a11050-9d9c-4487-a705-04263cfa461a.md", "e5a11050-9d9c-4487-a705-04263cfa461a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlMd2, [Type]::Missing, "e5a11050-9d9c-4487-a705-04263cfa461a.md", "e5a11050-9d9c-4487-a705-04263cfa461a.md")

$wsZh.Columns.Item(3).ColumnWidth = $wideColWidth
$wsZh.Columns.Item(9).ColumnWidth = $fullColWidth
$wsZh.Columns.Item(10).ColumnWidth = $fullColWidth

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("J2").Value = "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.252e131be62cc623e431e38256b7e8c7a708a19d.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-22 10:42:27"

$wsDe.Range("J3").Value = "e5a11050-9d9c-4487-a705-04263cfa461a.e97206c15f67c039fb036d86b15dd20e106ba7c3.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-22 10:42:27"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlMd1, [Type]::Missing, "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.md", "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlMd1, [Type]::Missing, "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.md", "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlMd2, [Type]::Missing, "e5a11050-9d9c-4487-a705-04263cfa461a.md", "e5a11050-9d9c-4487-a705-04263cfa461a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlMd2, [Type]::Missing, "e5a11050-9d9c-4487-a705-04263cfa461a.md", "e5a11050-9d9c-4487-a705-04263cfa461a.md")

$wsDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDe.Columns.Item(9).ColumnWidth = $fullColWidth
$wsDe.Columns.Item(10).ColumnWidth = $fullColWidth

# I2/I3 on both sheets now hold the "Latest Target File" md link (same file
# as column A for that row) -- set that last so the value write is not
# clobbered by Hyperlinks.Add (which already sets text/tooltip) and matches
# exactly what column A shows.
$wsZh.Range("I2").Value = "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.md"
$wsZh.Range("I3").Value = "e5a11050-9d9c-4487-a705-04263cfa461a.md"
$wsDe.Range("I2").Value = "24966354-00aa-4c6f-b8f1-c3a0a34cee9a.md"
$wsDe.Range("I3").Value = "e5a11050-9d9c-4487-a705-04263cfa461a.md"
